$wb = $excel.ActiveWorkbook

# The two sheets "展览" and "全部类型" contain identical data tables.
# Row 3 (F3): 想去人数 (interest count) 203 -> 205
# Row 5 (F5): 想去人数 (interest count) 51 -> 53
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 205
    $ws.Range("F5").Value = 53
}
